$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the existing
# header style (bold, bordered, centered) used by H1 etc.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I/J data columns, row by row.
$valuesI = @{ 2 = 5; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 8; 9 = 7; 10 = 6; 11 = 4 }
$valuesJ = @{ 2 = 8; 3 = 6; 4 = 5; 5 = 5; 6 = 6; 7 = 3; 8 = 9; 9 = 8; 10 = 7; 11 = 5 }

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 9).Value = $valuesI[$row]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$row]
}
